# Update "paises" workbook: refresh country case counts and the
# "Datos actualizados" timestamp (commit: "Update countries & provincias Spain")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (row 1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 01:22"

# --- Updated per-country rows -----------------------------------------------
# Each row keeps its position in the (descending, by "Casos totales") sorted
# table; when two countries' totals cross over, the rows swap which country
# they display along with its refreshed figures.

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 501648
$ws.Range("C4").Value = 33082
$ws.Range("D4").Value = 27239
$ws.Range("E4").Value = 455734
$ws.Range("F4").Value = 10916
$ws.Range("G4").Value = 1984
$ws.Range("H4").Value = 18675

# Row 16: Canada
$ws.Range("A16").Value = "Canada"
$ws.Range("B16").Value = 22148
$ws.Range("C16").Value = 1383
$ws.Range("D16").Value = 6013
$ws.Range("E16").Value = 15566
$ws.Range("F16").Value = 557
$ws.Range("G16").Value = 60
$ws.Range("H16").Value = 569

# Row 17: Brasil
$ws.Range("A17").Value = "Brasil"
$ws.Range("B17").Value = 19789
$ws.Range("C17").Value = 1644
$ws.Range("D17").Value = 173
$ws.Range("E17").Value = 18548
$ws.Range("F17").Value = 296
$ws.Range("G17").Value = 114
$ws.Range("H17").Value = 1068

# Row 91: now Costa de Marfil (overtakes Uruguay)
$ws.Range("A91").Value = "Costa de Marfil"
$ws.Range("B91").Value = 480
$ws.Range("C91").Value = 36
$ws.Range("D91").Value = 54
$ws.Range("E91").Value = 423
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 3

# Row 92: now Uruguay
$ws.Range("A92").Value = "Uruguay"
$ws.Range("B92").Value = 473
$ws.Range("C92").Value = 17
$ws.Range("D92").Value = 206
$ws.Range("E92").Value = 260
$ws.Range("F92").Value = 13
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 7

# Row 93: Burkina Faso
$ws.Range("A93").Value = "Burkina Faso"
$ws.Range("B93").Value = 448
$ws.Range("C93").Value = 5
$ws.Range("D93").Value = 149
$ws.Range("E93").Value = 273
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 2
$ws.Range("H93").Value = 26

# Row 110: now Montenegro (overtakes Vietnam)
$ws.Range("A110").Value = "Montenegro"
$ws.Range("B110").Value = 257
$ws.Range("C110").Value = 5
$ws.Range("D110").Value = 4
$ws.Range("E110").Value = 251
$ws.Range("F110").Value = 7
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 2

# Row 111: now Vietnam
$ws.Range("A111").Value = "Vietnam"
$ws.Range("B111").Value = 257
$ws.Range("C111").Value = 2
$ws.Range("D111").Value = 144
$ws.Range("E111").Value = 113
$ws.Range("F111").Value = 8
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 0

# Row 120: Venezuela
$ws.Range("A120").Value = "Venezuela"
$ws.Range("B120").Value = 175
$ws.Range("C120").Value = 4
$ws.Range("D120").Value = 84
$ws.Range("E120").Value = 82
$ws.Range("F120").Value = 6
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 9

# Row 165: now Antigua y Barbuda (overtakes Mozambique)
$ws.Range("A165").Value = "Antigua y Barbuda"
$ws.Range("B165").Value = 21
$ws.Range("C165").Value = 2
$ws.Range("D165").Value = 0
$ws.Range("E165").Value = 19
$ws.Range("F165").Value = 1
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 2

# Row 166: now Mozambique
$ws.Range("A166").Value = "Mozambique"
$ws.Range("B166").Value = 20
$ws.Range("C166").Value = 3
$ws.Range("D166").Value = 2
$ws.Range("E166").Value = 18
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 0
